$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "fas"
$ws.Range("B6").Value = "sff"
$ws.Range("B7").Value = "sf"
$ws.Range("B8").Value = 1
$ws.Range("B9").Value = 2
$ws.Range("B10").Value = 3
$ws.Range("B11").Value = 3
$ws.Range("B12").Value = "EUR"
$ws.Range("B13").Value = 4
$ws.Range("B14").Value = 6

$ws.Range("C15").Value = 0
$ws.Range("C16").Value = 0
$ws.Range("C17").Value = 0
$ws.Range("C18").Value = 0
$ws.Range("C19").Value = 0
$ws.Range("C20").Value = 0
$ws.Range("C21").Value = 0
$ws.Range("C22").Value = 0
$ws.Range("C23").Value = 0
$ws.Range("C24").Value = 0
$ws.Range("C25").Value = 0
